$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.824.04"
$ws.Range("E2").Value = "  +4.92%  "
$ws.Range("D3").Value = "1.611.49"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'213.67"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("E6").Value = "  +6.82%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").Value = "'26.87"
$ws.Range("E8").Value = "  +11.39%  "
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "1.842.44"
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("D13").Value = "1.608.81"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").Value = "29.850.17"
$ws.Range("E14").Value = "  +4.92%  "
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("E16").Value = "  +3.58%  "
$ws.Range("D17").Value = "'244.42"
$ws.Range("E17").Value = "  +6.71%  "
$ws.Range("D18").Value = "'63.41"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").Value = "'7.63"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "0.0₃0694"
$ws.Range("E20").Value = "  +3.25%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'4.05"
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("D23").Value = "'9.25"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("E24").Value = "  +3.91%  "
$ws.Range("D25").Value = "'156.00"
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("D26").Value = "'15.33"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("E27").Value = "  +5.31%  "
$ws.Range("D28").Value = "'6.40"
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'0.0473"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").Value = "1.442.25"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.84"
$ws.Range("E36").Value = "  +9.98%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "'0.0166"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("D41").Value = "'55.45"
$ws.Range("E41").Value = "  +29.16%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "'0.797"
$ws.Range("E43").Value = "  +3.19%  "
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "'0.0468"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").Value = "'65.92"
$ws.Range("E46").Value = "  +6.74%  "
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "1.753.81"
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("D49").Value = "'87.05"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  -4.15%  "
$ws.Range("D51").Value = "0.0₆0105"
$ws.Range("E51").Value = "  +3.48%  "
